# initial_request.docx template edit
#
# Three changes:
#  1. Split the "a {{ student.grade }}" sentence so that an
#     "After grade 12" case is special-cased with a new
#     {% if (student.grade != "After grade 12") %} ... block.
#  2. Add a new paragraph right after the "student in your district"
#     paragraph that closes the above if-block and adds an
#     {% if (student.grade == "After grade 12") %} branch.
#  3. Remove a stray <w:lastRenderedPageBreak/> before "{{ item".

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1
# ---------------------------------------------------------------
$rng1 = $d.Content
$ok1 = $rng1.Find.Execute("student.name }}{% endif %}, a {{ student.grade }} {% if (student.grade !=")
if (-not $ok1) {
    throw "Change 1 anchor text not found"
}
$target1 = $d.Range($rng1.Start, $rng1.End)

$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:r w:rsidRPr="0052748F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">student.name }}{% endif %}, </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{% if (</w:t></w:r>
<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0052748F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>student.grade</w:t></w:r><w:proofErr w:type="spellEnd"/>
<w:r w:rsidRPr="0052748F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> != &#8220;After grade 12&#8221;) %}</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">a {{ </w:t></w:r>
<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>student.grade</w:t></w:r><w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> }} {% if (</w:t></w:r>
<w:proofErr w:type="spellStart"/><w:r w:rsidRPr="0052748F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>student.grade</w:t></w:r><w:proofErr w:type="spellEnd"/>
<w:r w:rsidRPr="0052748F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> !=</w:t></w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target1.InsertXML($xml1)

# ---------------------------------------------------------------
# Change 2 - insert a brand-new paragraph right after the paragraph
# that ends in "...student in your district{% endif %}."
# ---------------------------------------------------------------
$rng2 = $d.Content
$ok2 = $rng2.Find.Execute("% else %}student in your district{% endif %}.")
if (-not $ok2) {
    throw "Change 2 anchor text not found"
}
$insertPoint2 = $d.Range($rng2.End, $rng2.End)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{% endif %}</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{% if (</w:t></w:r>
<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>student.grade</w:t></w:r><w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> == &#8220;After grade 12&#8221;) %}an after grade 12 student at your school</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{% endif %}</w:t></w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertPoint2.InsertXML($xml2)

# ---------------------------------------------------------------
# Change 3 - drop the stray <w:lastRenderedPageBreak/> in front of
# "{{ item" (spans the "{{ item" and " }}" runs plus the proofErr
# sitting between them so nothing shifts out of place).
# ---------------------------------------------------------------
$rng3 = $d.Content
$ok3 = $rng3.Find.Execute("{{ item }}")
if (-not $ok3) {
    throw "Change 3 anchor text not found"
}
$target3 = $d.Range($rng3.Start, $rng3.End)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:r w:rsidRPr="009F027E"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{{ item</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r w:rsidRPr="009F027E"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target3.InsertXML($xml3)

Write-Output "Done."
